# Auto-generated edit script: correct 26 grammatical-error strings in column E
# of Sheet1, then restore the sheet view's scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 'There is no explicit definition of the strategy in the teams. The team does not have a shared accessible document where they can see the strategy. Alignement is purely coincidental and not deliberate'
$ws.Range("E5").Value = 'There is a review of alignment to strategy in the Post implementation PI sessions but no review during the implementation stage.'
$ws.Range("E14").Value = 'There is visible innovation/optimisation in the team in line with the strategy. These are captured as cards on the team''s story wall and as extra achievements in the scope management tool.'
$ws.Range("E36").Value = 'The retrospectives are run every iteration and the stakeholders are present in the retrospectives. Stakeholders contribute to retrospectives and action items, which are allocated to them are tracked by the team on the team wall.'
$ws.Range("E53").Value = 'There are continuous assessments of the code metrics to ensure the team is aligned to the quality standards. Backlog items are created to address gaps identified from the code metrics'' perspective and these are tracked and prioritized.'
$ws.Range("E54").Value = 'Technical debt items are tracked and the code static analysis metrics are tracked to show progress with respect to technical debt reduction.'
$ws.Range("E55").Value = 'Artefacts, code and reports are kept in developer machines. There is no use of a version management technology.'
$ws.Range("E60").Value = 'The environments required by the development teams (developers, QAs, UAT etc) are provisioned manually and shared by multiple teams. The team is unable to modify configuration which they want to test ad hoc and needs to raise change requests.'
$ws.Range("E61").Value = 'The source code is checked in to a version control management tool even though the team commits are very infrequent.'
$ws.Range("E68").Value = 'Build statuses are shown on dashboards which the team can see. Team has a view of the broken builds history. Team continues with coding tasks even when the build is broken. The team needs to start moving towards a state where a broken build leads to a stop in new functionality until the build is fixed (Think Broken Window Syndrome).'
$ws.Range("E91").Value = 'There is very minimal handoffs. The resolution is effective and involves very few teams, generally the core team and the dependent system team.'
$ws.Range("E108").Value = 'The risks are allocated owners as well as mitigation/action plans. These are clearly captured in the story management system.'
$ws.Range("E114").Value = 'Roles, responsibilities and accountabilities for the risks owners are defined, agreed upon and documented where it is accessible to all team members.'
$ws.Range("E122").Value = 'At this level of proficiency, the team does a lot of Big Up-Front designs for the systems as opposed to a "just enough" design approach. There are several ceremonial design sessions for the entire application without necessarily taking into effect the fact that there might need to be changes should the product vision change.'
$ws.Range("E126").Value = 'Design issues identified are tracked and prioritised in the backlog and attended to, similar to the technical debts.'
$ws.Range("E129").Value = 'The technical design is owned by the feature teams. The technical architect is part of the core team and is not a stakeholder, external to the team and making decisions on their behalf.'
$ws.Range("E138").Value = 'The teams work on functionality without necessarily knowing what the goal and vision is. Some people in the team or some stakeholders may know the vision but this is not shared across and this has not been documented. The goals/vision need to be documented and easily accessible.'
$ws.Range("E156").Value = 'The release dates are seldom moved out. The team is able to reduce the release features rather than the move out the release dates.'
$ws.Range("E157").Value = 'There is a well defined process for handling the release of the artefacts to operations and support teams in place. This is documented and well understood by all the team members. The team at this level is investigating means of embedding operations/support in the team structure.'
$ws.Range("E211").Value = 'Database migrations happen as part of deployment. The deployment pipelines have steps to ensure that the database base is in a state expected by the application. These scripts are idempotent in nature.'
$ws.Range("E219").Value = 'Multiple test environments are readily available for the exclusive use of the team, including a production-like environment that allows a reasonable level of non-functional and systems level testing. These environments are controlled by the teams.'
$ws.Range("E222").Value = 'Data and environments are versioned and managed in the version control management system in the same way as source code and artefacts.'
$ws.Range("E231").Value = 'At this level the team is not cross-functional in nature. Some of the skills required to complete end-to-end customer features are missing and there is dependence on teams/people outside the team.'
$ws.Range("E239").Value = 'Cross-functional and cross-component skills are embedded in the teams and all the requisite skills required to deliver end-to-end customer value are embedded in the teams. There is no external team dependency required to deliver the goals and vision of the product.'
$ws.Range("E241").Value = 'The triad (business, technology and quality) is fully functional in the teams and the product owners within the teams are empowered to make the decisions pertaining to the product roadmap without the need of external consultation or validation.'
$ws.Range("E244").Value = 'Team has autonomy and is self-organized. The team composition is well understood by all team members and stakeholders, and the roles and responsibilities are clearly defined and understood.'

# Restore the window's scroll position (was showing row 242 at top,
# now showing row 238 at top) while keeping the existing selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 238
$win.ScrollColumn = 1

